$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 494 (sheet is sorted alphabetically by column A;
# "The Civil Conflict Ceasefire dataset" sorts between
# "The CIRI Human Rights Dataset" and "The Comparative Legislators Database").
$ws.Rows(494).Insert()

# Fill in the new row's data.
$ws.Range("A494").Value = "The Civil Conflict Ceasefire dataset"
$ws.Range("B494").Value = "international relations"
$ws.Range("C494").Value = "https://ceasefireproject.org/download/"
$ws.Range("D494").Value = "ceasefire, civil conflict, conflict management"

$ws.Range("F494").Value = 1
$ws.Range("G494").Value = 1
$ws.Range("H494").Value = 1
$ws.Range("I494").Value = 1
$ws.Range("J494").Value = 1

$ws.Range("K494").Value = 1989
$ws.Range("L494").Value = 2020

$ws.Range("M494").Value = "online"
$ws.Range("N494").Value = "no"
$ws.Range("O494").Value = 1

$ws.Range("P494").Value = "https://ceasefireproject.org/wp-content/uploads/2022/11/Civil-Conflict-Ceasefire-Data-set-%E2%80%93-Codebook-V1.pdf"
$ws.Range("R494").Value = "https://ceasefireproject.org/wp-content/uploads/2022/11/CFD_oct_2022_id-1.dta"
$ws.Range("T494").Value = "https://ceasefireproject.org/wp-content/uploads/2022/11/CFD_oct_2022_id-1.xlsx"

$ws.Range("W494").Value = "location"
$ws.Range("X494").Value = "end_yr"
$ws.Range("Z494").Value = "10.1177/00220027221129183"

$ws.Range("AB494").Value = 20230909

# Hyperlinks for the new row (matches order of hyperlinks added in the source edit).
$ws.Hyperlinks.Add($ws.Range("C494"), $ws.Range("C494").Value2)
$ws.Hyperlinks.Add($ws.Range("P494"), $ws.Range("P494").Value2)
$ws.Hyperlinks.Add($ws.Range("T494"), $ws.Range("T494").Value2)
$ws.Hyperlinks.Add($ws.Range("R494"), $ws.Range("R494").Value2)
